{"js": "// Update the date line and all twenty-five division problems in the table\n// to the \"next day\" worksheet values, per the commit's regenerated output.\n\nconst replacements = [\n  [\"2024-11-22 Friday\", \"2024-11-23 Saturday\"],\n  [\"114\u00f79=\", \"348\u00f77=\"],\n  [\"965\u00f77=\", \"958\u00f78=\"],\n  [\"135\u00f75=\", \"948\u00f74=\"],\n  [\"231\u00f77=\", \"300\u00f72=\"],\n  [\"325\u00f73=\", \"675\u00f73=\"],\n  [\"855\u00f75=\", \"914\u00f74=\"],\n  [\"470\u00f73=\", \"740\u00f77=\"],\n  [\"397\u00f72=\", \"617\u00f75=\"],\n  [\"661\u00f78=\", \"769\u00f79=\"],\n  [\"102\u00f76=\", \"114\u00f72=\"],\n  [\"106\u00f79=\", \"511\u00f73=\"],\n  [\"988\u00f74=\", \"453\u00f79=\"],\n  [\"681\u00f77=\", \"956\u00f79=\"],\n  [\"196\u00f75=\", \"748\u00f75=\"],\n  [\"478\u00f79=\", \"462\u00f73=\"],\n  [\"177\u00f75=\", \"450\u00f77=\"],\n  [\"686\u00f79=\", \"358\u00f73=\"],\n  [\"592\u00f72=\", \"590\u00f76=\"],\n  [\"626\u00f72=\", \"872\u00f77=\"],\n  [\"288\u00f73=\", \"615\u00f74=\"],\n  [\"373\u00f75=\", \"499\u00f72=\"],\n  [\"899\u00f78=\", \"502\u00f77=\"],\n  [\"271\u00f73=\", \"949\u00f72=\"],\n  [\"955\u00f74=\", \"362\u00f77=\"],\n  [\"987\u00f78=\", \"664\u00f72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and all twenty-five division problems in the table\n# to the \"next day\" worksheet values, per the commit's regenerated output.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-11-22 Friday\", \"2024-11-23 Saturday\"),\n    @(\"114\u00f79=\", \"348\u00f77=\"),\n    @(\"965\u00f77=\", \"958\u00f78=\"),\n    @(\"135\u00f75=\", \"948\u00f74=\"),\n    @(\"231\u00f77=\", \"300\u00f72=\"),\n    @(\"325\u00f73=\", \"675\u00f73=\"),\n    @(\"855\u00f75=\", \"914\u00f74=\"),\n    @(\"470\u00f73=\", \"740\u00f77=\"),\n    @(\"397\u00f72=\", \"617\u00f75=\"),\n    @(\"661\u00f78=\", \"769\u00f79=\"),\n    @(\"102\u00f76=\", \"114\u00f72=\"),\n    @(\"106\u00f79=\", \"511\u00f73=\"),\n    @(\"988\u00f74=\", \"453\u00f79=\"),\n    @(\"681\u00f77=\", \"956\u00f79=\"),\n    @(\"196\u00f75=\", \"748\u00f75=\"),\n    @(\"478\u00f79=\", \"462\u00f73=\"),\n    @(\"177\u00f75=\", \"450\u00f77=\"),\n    @(\"686\u00f79=\", \"358\u00f73=\"),\n    @(\"592\u00f72=\", \"590\u00f76=\"),\n    @(\"626\u00f72=\", \"872\u00f77=\"),\n    @(\"288\u00f73=\", \"615\u00f74=\"),\n    @(\"373\u00f75=\", \"499\u00f72=\"),\n    @(\"899\u00f78=\", \"502\u00f77=\"),\n    @(\"271\u00f73=\", \"949\u00f72=\"),\n    @(\"955\u00f74=\", \"362\u00f77=\"),\n    @(\"987\u00f78=\", \"664\u00f72=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
